$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (matching the workbook's original
# inline-string cells) without leaving a stray "@" number-format style
# behind on the cell - Excel auto-detects numeric-looking strings (e.g.
# "218.51") and would otherwise store them as real numbers.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Rows where only the Price (D) and/or Volume(1h) (E) columns changed.
# $null means "leave that column unchanged".
$rowUpdates = @{
    2  = @{ D = "26.744.74";      E = "  -6.75%  " }
    3  = @{ D = "1.696.09";       E = "  -5.99%  " }
    4  = @{ D = $null;            E = "  +0.36%  " }
    5  = @{ D = "218.51";         E = "  -5.49%  " }
    6  = @{ D = "0.5037";         E = "  -15.58%  " }
    7  = @{ D = $null;            E = "  +0.25%  " }
    8  = @{ D = "0.2591";         E = "  -6.90%  " }
    9  = @{ D = "21.77";          E = "  -6.97%  " }
    10 = @{ D = "0.06108";        E = "  -10.78%  " }
    11 = @{ D = "0.07310";        E = "  -2.84%  " }
    12 = @{ D = "1.693.99";       E = "  -5.92%  " }
    13 = @{ D = "4.417";          E = "  -6.16%  " }
    14 = @{ D = "1.927.96";       E = "  -5.91%  " }
    15 = @{ D = "0.5718";         E = "  -8.77%  " }
    16 = @{ D = "0.000008154";    E = "  -11.41%  " }
    17 = @{ D = "65.02";          E = "  -13.68%  " }
    18 = @{ D = "26.779.01";      E = "  -6.50%  " }
    19 = @{ D = "5.007";          E = "  -8.44%  " }
    20 = @{ D = "1.006";          E = "  +0.31%  " }
    21 = @{ D = "10.74";          E = "  -6.16%  " }
    22 = @{ D = "183.78";         E = "  -12.75%  " }
    23 = @{ D = "6.207";          E = "  -9.33%  " }
    24 = @{ D = $null;            E = "  +0.36%  " }
    25 = @{ D = "145.16";         E = "  -6.06%  " }
    26 = @{ D = "7.613";          E = "  -2.92%  " }
    27 = @{ D = "0.1140";         E = "  -10.80%  " }
    28 = @{ D = "15.24";          E = "  -6.92%  " }
    29 = @{ D = "1.316";          E = "  -9.05%  " }
    30 = @{ D = "0.05602";        E = "  -10.36%  " }
    31 = @{ D = "1.328";          E = "  -6.49%  " }
    32 = @{ D = "3.459";          E = "  -8.02%  " }
    33 = @{ D = "3.434";          E = "  -7.97%  " }
    34 = @{ D = "1.655";          E = "  -3.66%  " }
    35 = @{ D = $null;            E = "  -4.63%  " }
    36 = @{ D = "2.408";          E = "  -3.49%  " }
    37 = @{ D = "0.5871";         E = "  -7.87%  " }
    38 = @{ D = "2.634";          E = "  -3.09%  " }
    39 = @{ D = "0.01582";        E = "  -7.36%  " }
    40 = @{ D = "1.069.84";       E = $null }
    41 = @{ D = "5.886";          E = "  -7.97%  " }
    42 = @{ D = "0.8513";         E = "  -1.72%  " }
    43 = @{ D = $null;            E = "  +0.12%  " }
    44 = @{ D = "98.46";          E = "  -2.17%  " }
    45 = @{ D = "1.857.16";       E = "  -5.30%  " }
    46 = @{ D = "56.30";          E = "  -7.08%  " }
    49 = @{ D = "8.052";          E = "  -3.39%  " }
    50 = @{ D = "0.4336";         E = "  -3.66%  " }
    51 = @{ D = "0.05205";        E = "  -4.32%  " }
}

foreach ($row in $rowUpdates.Keys) {
    $vals = $rowUpdates[$row]
    if ($null -ne $vals.D) {
        Set-TextValue $row 4 $vals.D
    }
    if ($null -ne $vals.E) {
        Set-TextValue $row 5 $vals.E
    }
}

# Rows 47 and 48 swap their coin content entirely: the former Frax row (47)
# becomes BabyDogeCoin, and the former BabyDogeCoin row (48) becomes Frax,
# each with refreshed price/volume figures.
Set-TextValue 47 2 "BabyDogeCoin"
Set-TextValue 47 3 "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue 47 4 "0.00000000106"
Set-TextValue 47 5 "  -5.56%  "

Set-TextValue 48 2 "Frax"
Set-TextValue 48 3 "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue 48 4 "1.003"
Set-TextValue 48 5 "  -0.26%  "
